$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Genome)
$ws.Range("E3").Value = 67.32348111658456
$ws.Range("F3").Value = 74.13793103448276
$ws.Range("G3").Value = 72.41379310344827
$ws.Range("H3").Value = 84.40065681444992
$ws.Range("I3").Value = 87.25451559934319

# Row 4 (Realm)
$ws.Range("E4").Value = 75.94827586206897
$ws.Range("F4").Value = 80.94827586206897
$ws.Range("G4").Value = 81.37931034482759
$ws.Range("H4").Value = 88.70689655172413
$ws.Range("I4").Value = 92.57413793103451

# Row 5 (Kingdom)
$ws.Range("E5").Value = 73.48877374784111
$ws.Range("F5").Value = 78.75647668393782
$ws.Range("G5").Value = 78.41105354058722
$ws.Range("H5").Value = 85.49222797927462
$ws.Range("I5").Value = 90.95682210708118

# Row 6 (Phylum)
$ws.Range("E6").Value = 61.5916955017301
$ws.Range("F6").Value = 56.74740484429066
$ws.Range("G6").Value = 55.88235294117647
$ws.Range("H6").Value = 71.280276816609
$ws.Range("I6").Value = 83.4134948096886

# Row 7 (Class)
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 5845
$ws.Range("E7").Value = 51.15483319076134
$ws.Range("F7").Value = 52.95124037639007
$ws.Range("G7").Value = 47.56201881950385
$ws.Range("H7").Value = 63.47305389221557
$ws.Range("I7").Value = 80.22754491017963

# Row 8 (Order)
$ws.Range("C8").Value = 48
$ws.Range("D8").Value = 5838
$ws.Range("E8").Value = 48.88698630136986
$ws.Range("F8").Value = 55.65068493150685
$ws.Range("G8").Value = 48.88698630136986
$ws.Range("H8").Value = 60.61643835616438
$ws.Range("I8").Value = 79.62157534246576

# Row 9 (Family)
$ws.Range("C9").Value = 102
$ws.Range("D9").Value = 5990
$ws.Range("E9").Value = 36.64440734557596
$ws.Range("F9").Value = 43.23873121869783
$ws.Range("G9").Value = 27.04507512520868
$ws.Range("H9").Value = 42.98831385642738
$ws.Range("I9").Value = 74.45909849749583

# Row 10 (Genus)
$ws.Range("C10").Value = 360
$ws.Range("D10").Value = 4673
$ws.Range("E10").Value = 44.59893048128342
$ws.Range("F10").Value = 36.79144385026738
$ws.Range("G10").Value = 18.82352941176471
$ws.Range("H10").Value = 17.64705882352941
$ws.Range("I10").Value = 68.71229946524065
